# jo/update working example full
#
# The "Z1" row (row 8: tag=Z1, feature_name=Z1, init_values=120) is removed
# from the InitialValues sheet. Deleting the entire row shifts every row
# below it up by one, which is exactly what the target diff shows (the now
# out-of-use "Z1" shared string is dropped, dimension shrinks from C18 to
# C17, and all the tag/feature_name/init_values rows below slide up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire "Z1" row (row 8) — pulls every following row up by one.
$ws.Rows.Item(8).Delete()

# Match the post-edit selection state recorded in the workbook (cursor left
# on H13 after the row delete).
$ws.Range("H13").Select()
